# Apply "Аналитик" dataset update to report.xlsx
# Sheet 1: "Статистика по годам" (Stats by year)
# Sheet 2: "Статистика по городам" (Stats by city)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Статистика по годам")
$ws2 = $wb.Worksheets.Item("Статистика по городам")

# ---------------------------------------------------------------------------
# Sheet 1 - header label updates (React -> Аналитик)
# ---------------------------------------------------------------------------
$ws1.Range("C1").Value = "Средняя зарплата - Аналитик"
$ws1.Range("E1").Value = "Количество вакансий - Аналитик"

# ---------------------------------------------------------------------------
# Sheet 1 - grow the data block from 3 rows (A2:E4) to 16 rows (A2:E17) and
# replace all of the year/salary/vacancy figures.
# First propagate the existing bordered-cell style (style of row 4) down over
# the whole new range, then write the values on top of it.
# ---------------------------------------------------------------------------
$ws1.Range("A4:E4").Copy()
$ws1.Range("A5:E17").PasteSpecial(-4122)   # xlPasteFormats

$sheet1Data = @(
  @(2007, 38916, 40641,   2196,   34),
  @(2008, 43646, 48428,  17549,  196),
  @(2009, 42492, 48109,  17709,  171),
  @(2010, 43846, 49577,  29093,  328),
  @(2011, 47451, 52794,  36700,  418),
  @(2012, 48243, 58341,  44153,  374),
  @(2013, 51510, 57004,  59954,  420),
  @(2014, 50658, 58768,  66837,  504),
  @(2015, 52696, 53326,  70039,  749),
  @(2016, 62675, 54236,  75145,  911),
  @(2017, 60935, 56558,  82823, 1201),
  @(2018, 58335, 61080, 131701, 1578),
  @(2019, 69467, 71288, 115086, 1482),
  @(2020, 73431, 80145, 102243, 1349),
  @(2021, 82690, 87473,  57623,  805),
  @(2022, 91795, 91340,  18294,  305)
)

$rowCount = $sheet1Data.Count
$colCount = 5
$arr1 = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
  for ($j = 0; $j -lt $colCount; $j++) {
    $arr1[$i,$j] = $sheet1Data[$i][$j]
  }
}
$ws1.Range("A2:E17").Value = $arr1

# ---------------------------------------------------------------------------
# Sheet 1 - column width changes (C and E got a bit wider)
# ColumnWidth setter adds ~0.8333 "padding" before it is stored back to the
# OOXML width attribute, so back that constant out to land on the exact
# target width.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 29 - 0.8333333333333333
$ws1.Columns.Item(5).ColumnWidth = 32 - 0.8333333333333333

# ---------------------------------------------------------------------------
# Sheet 2 - header row tweak: the spacer column (C) now just holds a single
# space like the rest of the rows below it, and the "Город"/"Доля вакансий"
# labels shift one column to the right (D/E) to line up with the data.
# ---------------------------------------------------------------------------
$ws2.Range("C1").Value = " "
$ws2.Range("D1").Value = "Город"
$ws2.Range("E1").Value = "Доля вакансий"

# ---------------------------------------------------------------------------
# Sheet 2 - refreshed city figures (left block: city/salary, right block:
# city/vacancy share). Rows are re-sorted by the new numbers.
# ---------------------------------------------------------------------------
$sheet2Left = @(
  @("Москва",          76970),
  @("Санкт-Петербург", 65286),
  @("Новосибирск",     62254),
  @("Екатеринбург",    60962),
  @("Казань",          52580),
  @("Краснодар",       51644),
  @("Челябинск",       51265),
  @("Самара",          50994),
  @("Пермь",           48089),
  @("Нижний Новгород", 47662)
)

$sheet2Right = @(
  @("Москва",          0.3246),
  @("Санкт-Петербург", 0.1197),
  @("Новосибирск",     0.0271),
  @("Казань",          0.0237),
  @("Нижний Новгород", 0.0232),
  @("Ростов-на-Дону",  0.0209),
  @("Екатеринбург",    0.0207),
  @("Краснодар",       0.0185),
  @("Самара",          0.0143),
  @("Воронеж",         0.0141)
)

$n = $sheet2Left.Count
$arrLeft = New-Object 'object[,]' $n,2
$arrRight = New-Object 'object[,]' $n,2
for ($i = 0; $i -lt $n; $i++) {
  $arrLeft[$i,0]  = $sheet2Left[$i][0]
  $arrLeft[$i,1]  = $sheet2Left[$i][1]
  $arrRight[$i,0] = $sheet2Right[$i][0]
  $arrRight[$i,1] = $sheet2Right[$i][1]
}
$ws2.Range("A2:B11").Value = $arrLeft
$ws2.Range("D2:E11").Value = $arrRight

# ---------------------------------------------------------------------------
# Sheet 2 - column width changes (spacer column C narrower, share column E
# wider)
# ---------------------------------------------------------------------------
$ws2.Columns.Item(3).ColumnWidth = 2 - 0.8333333333333333
$ws2.Columns.Item(5).ColumnWidth = 15 - 0.8333333333333333
